# Apply weekly update to the "Hortaliza, Macroferia Regional de Talca - Brocoli" sheet.
# A brand-new weekly record is inserted as row 167 (pushing all the following
# records down by one row); the new record repeats the values that used to
# sit in the former row 167, except for a newer report date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 167

# Insert a new blank row at position 167; Excel shifts rows 167..200 down to 168..201
# and the sheet's used range / dimension grows to A1:R201 automatically.
$ws.Rows.Item($insertRow).Insert()

# Populate the freshly inserted row with the same data the old row 167 held,
# only the "Fecha" (date) column changes to the new reporting date.
$ws.Cells.Item($insertRow, 1).Value  = 5
$ws.Cells.Item($insertRow, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($insertRow, 3).Value  = "Maule"
$ws.Cells.Item($insertRow, 4).Value  = "10/07/2021"
$ws.Cells.Item($insertRow, 5).Value  = 7
$ws.Cells.Item($insertRow, 6).Value  = 100112023
$ws.Cells.Item($insertRow, 7).Value  = "Brócoli"
$ws.Cells.Item($insertRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($insertRow, 9).Value  = "Primera"
$ws.Cells.Item($insertRow, 10).Value = 3000
$ws.Cells.Item($insertRow, 11).Value = 600
$ws.Cells.Item($insertRow, 12).Value = 600
$ws.Cells.Item($insertRow, 13).Value = 600
$ws.Cells.Item($insertRow, 14).Value = "`$/unidad"
$ws.Cells.Item($insertRow, 15).Value = "Región del Maule"
$ws.Cells.Item($insertRow, 16).Value = 600
$ws.Cells.Item($insertRow, 17).Value = 1
$ws.Cells.Item($insertRow, 18).Value = "Hortaliza"
